$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'schubert-winterreise_145'
$ws.Range("B2").Value = 'schubert-winterreise_61'
$ws.Range("C2").Value = 0.323076923076923
$ws.Range("D2").Value = '[[''D:maj/A'', ''G:maj'', ''D:maj/A'']]'
$ws.Range("E2").Value = '[[''G:maj'', ''C:maj/G'', ''G:maj'']]'
$ws.Range("F2").Value = '[(143.58, 148.16)]'
$ws.Range("G2").Value = '[(18.12, 24.54)]'
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'

# Row 3
$ws.Range("A3").Value = 'schubert-winterreise_75'
$ws.Range("B3").Value = 'schubert-winterreise_94'
$ws.Range("C3").Value = 0.1098901098901099
$ws.Range("D3").Value = '[[''C:7'', ''F:min'', ''A#:min/C#''], [''F:min/C'', ''C'', ''F:min/C'']]'
$ws.Range("E3").Value = '[[''F:7/A#'', ''A#:min'', ''D#:min/A#''], [''A#:min'', ''F:maj'', ''A#:min'']]'
$ws.Range("F3").Value = '[(24.42, 32.62), (45.72, 49.9)]'
$ws.Range("G3").Value = '[(27.52, 35.76), (15.6, 24.08)]'
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'

# Row 4
$ws.Range("A4").Value = 'isophonics_270'
$ws.Range("B4").Value = 'isophonics_150'
$ws.Range("C4").Value = 0.1052631578947368
$ws.Range("D4").Value = '[[''G'', ''E:min'', ''A'']]'
$ws.Range("E4").Value = '[[''G'', ''E:min'', ''A'']]'
$ws.Range("F4").Value = '[(5.038752, 9.647913)]'
$ws.Range("G4").Value = '[(50.135784, 55.151295)]'
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""

# Row 5
$ws.Range("A5").Value = 'isophonics_132'
$ws.Range("B5").Value = 'isophonics_295'
$ws.Range("C5").Value = 0.1174242424242424
$ws.Range("D5").Value = '[[''F#'', ''B'', ''B/7'']]'
$ws.Range("E5").Value = '[[''G'', ''C/5'', ''C'']]'
$ws.Range("F5").Value = '[(37.936167, 41.77907)]'
$ws.Range("G5").Value = '[(20.870746, 26.837029)]'
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""

# Row 6
$ws.Range("A6").Value = 'isophonics_21'
$ws.Range("B6").Value = 'isophonics_157'
$ws.Range("C6").Value = 0.3170289855072464
$ws.Range("D6").Value = '[[''C'', ''G/3'', ''C'', ''G/3'', ''C''], [''F'', ''F:min'', ''C'', ''G/3'', ''C'']]'
$ws.Range("E6").Value = '[[''E'', ''B'', ''E'', ''B'', ''E''], [''A'', ''A:min/b3'', ''E'', ''B'', ''E'']]'
$ws.Range("F6").Value = '[(130.148, 132.075), (38.588, 46.22)]'
$ws.Range("G6").Value = '[(89.257527, 110.178707), (23.370907, 35.305963)]'
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""

# Row 7
$ws.Range("A7").Value = 'schubert-winterreise_129'
$ws.Range("B7").Value = 'schubert-winterreise_191'
$ws.Range("C7").Value = 0.2363636363636364
$ws.Range("D7").Value = '[[''F:min'', ''C:7'', ''F:min''], [''F:min'', ''C:maj'', ''F:min'']]'
$ws.Range("E7").Value = '[[''B:min'', ''F#:7/C#'', ''B:min/D''], [''B:min'', ''F#:maj/A#'', ''B:min'']]'
$ws.Range("F7").Value = '[(8.2, 12.66), (0.78, 5.26)]'
$ws.Range("G7").Value = '[(0.68, 2.98), (16.26, 19.34)]'
$ws.Range("H7").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Range("I7").Value = ""

# Row 8
$ws.Range("A8").Value = 'isophonics_223'
$ws.Range("B8").Value = 'isophonics_99'
$ws.Range("C8").Value = 0.09519230769230769
$ws.Range("D8").Value = '[[''A'', ''F#:min'', ''D'']]'
$ws.Range("E8").Value = '[[''E'', ''C#:min'', ''A'']]'
$ws.Range("F8").Value = '[(25.908231, 32.885827)]'
$ws.Range("G8").Value = '[(23.574625, 26.08238)]'
$ws.Range("H8").Value = 'spotify:track:3KfbEIOC7YIv90FIfNSZpo'
$ws.Range("I8").Value = ""

# Row 9
$ws.Range("A9").Value = 'schubert-winterreise_86'
$ws.Range("B9").Value = 'isophonics_152'
$ws.Range("C9").Value = 0.1940559440559441
$ws.Range("D9").Value = '[[''D:7'', ''G:maj'', ''G:maj'']]'
$ws.Range("E9").Value = '[[''B:7'', ''E'', ''A/3'']]'
$ws.Range("F9").Value = '[(5.88, 11.46)]'
$ws.Range("G9").Value = '[(45.712, 51.648)]'
$ws.Range("H9").Value = 'spotify:track:0XfunCHFEeQnzm4NaY8rJr'
$ws.Range("I9").Value = ""

# Row 10
$ws.Range("A10").Value = 'schubert-winterreise_92'
$ws.Range("B10").Value = 'schubert-winterreise_62'
$ws.Range("C10").Value = 0.2363636363636364
$ws.Range("D10").Value = '[[''B:min'', ''F#:7/C#'', ''B:min/D'']]'
$ws.Range("E10").Value = '[[''A:min'', ''E:7/G#'', ''A:min'']]'
$ws.Range("F10").Value = '[(0.36, 2.44)]'
$ws.Range("G10").Value = '[(40.12, 47.88)]'
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = 'spotify:track:1yerCi2iQCVkdHG6rdRn7R'

# Row 11
$ws.Range("A11").Value = 'schubert-winterreise_170'
$ws.Range("B11").Value = 'schubert-winterreise_36'
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = '[[''B:7'', ''A:min/C'', ''D:7'', ''G:maj'', ''G:maj'', ''E:min'', ''A:min/C'', ''D:7'', ''G:maj'', ''C:maj/G'', ''G:maj'', ''A:min/C'', ''D:7'', ''G:maj'', ''A:min/C'', ''D:7'', ''G:maj'', ''D:min/F'', ''E:hdim7'', ''D:min/F'', ''E:(3,b5,b7)'', ''D:min/F'', ''D:min/F'', ''G:7'', ''C:min'', ''D:(3,b5,b7)'', ''C:min/D#'', ''F:min/G#'', ''G:hdim7'', ''F:min/G#'', ''G:(3,b5,b7)'', ''F:min/G#'', ''C:min/D#'', ''G:min/D'', ''D:7'', ''G:min'', ''A:(3,b5,b7)/G'', ''G:min'', ''G#:maj/C'', ''G:min'', ''A:(3,b5,b7)/G'', ''G:min'', ''D:7'', ''G:min'', ''G:maj'', ''G:7/F'', ''C:maj/E'', ''G:maj/D'', ''G:7/F'', ''C:maj/E'', ''G:maj/D'', ''A:min7/C'', ''D:7'', ''G:maj'', ''D:7/C'', ''G:maj/B'', ''D:7/C'', ''G:maj/B'', ''C:min'', ''G:min/A#'', ''C:min'', ''G:min/A#'', ''A:hdim7/D#'', ''D:7'', ''G:min'']]'
$ws.Range("E11").Value = '[[''B:7'', ''A:min/C'', ''D:7'', ''G:maj'', ''G:maj'', ''E:min'', ''A:min/C'', ''D:7'', ''G:maj'', ''C:maj/G'', ''G:maj'', ''A:min/C'', ''D:7'', ''G:maj'', ''A:min/C'', ''D:7'', ''G:maj'', ''D:min/F'', ''E:hdim7'', ''D:min/F'', ''E:(3,b5,b7)'', ''D:min/F'', ''D:min/F'', ''G:7'', ''C:min'', ''D:(3,b5,b7)'', ''C:min/D#'', ''F:min/G#'', ''G:hdim7'', ''F:min/G#'', ''G:(3,b5,b7)'', ''F:min/G#'', ''C:min/D#'', ''G:min/D'', ''D:7'', ''G:min'', ''A:(3,b5,b7)/G'', ''G:min'', ''G#:maj/C'', ''G:min'', ''A:(3,b5,b7)/G'', ''G:min'', ''D:7'', ''G:min'', ''G:maj'', ''G:7/F'', ''C:maj/E'', ''G:maj/D'', ''G:7/F'', ''C:maj/E'', ''G:maj/D'', ''A:min7/C'', ''D:7'', ''G:maj'', ''D:7/C'', ''G:maj/B'', ''D:7/C'', ''G:maj/B'', ''C:min'', ''G:min/A#'', ''C:min'', ''G:min/A#'', ''A:hdim7/D#'', ''D:7'', ''G:min'']]'
$ws.Range("F11").Value = '[(4.76, 118.78)]'
$ws.Range("G11").Value = '[(7.54, 122.26)]'
$ws.Range("H11").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Range("I11").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'

# Row 12
$ws.Range("A12").Value = 'isophonics_69'
$ws.Range("B12").Value = 'isophonics_107'
$ws.Range("C12").Value = 0.3176470588235294
$ws.Range("D12").Value = '[[''E'', ''E'', ''A'', ''E'']]'
$ws.Range("E12").Value = '[[''E'', ''E'', ''A'', ''E'']]'
$ws.Range("F12").Value = '[(22.125076, 36.625937)]'
$ws.Range("G12").Value = '[(122.976598, 134.098957)]'
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""

# Row 13
$ws.Range("A13").Value = 'isophonics_157'
$ws.Range("B13").Value = 'isophonics_1'
$ws.Range("C13").Value = 0.2724358974358975
$ws.Range("D13").Value = '[[''A'', ''E'', ''A'', ''E'', ''B'']]'
$ws.Range("E13").Value = '[[''Ab/5'', ''Eb'', ''Ab/5'', ''Eb'', ''Bb/3'']]'
$ws.Range("F13").Value = '[(68.150589, 83.20873)]'
$ws.Range("G13").Value = '[(18.978, 26.66)]'
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = ""

# Row 14
$ws.Range("A14").Value = 'schubert-winterreise_6'
$ws.Range("B14").Value = 'jaah_85'
$ws.Range("C14").Value = 0.08684863523573201
$ws.Range("D14").Value = '[[''B:7'', ''E:maj/B'', ''E:min/B'']]'
$ws.Range("E14").Value = '[[''Ab:7'', ''Db'', ''Db:min'']]'
$ws.Range("F14").Value = '[(4.54, 12.1)]'
$ws.Range("G14").Value = '[(36.66, 38.43)]'
$ws.Range("H14").Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'
$ws.Range("I14").Value = ""

# Row 15
$ws.Range("A15").Value = 'schubert-winterreise_157'
$ws.Range("B15").Value = 'jaah_85'
$ws.Range("C15").Value = 0.1136044880785414
$ws.Range("D15").Value = '[[''C:7'', ''F:maj'', ''F:maj/A''], [''F:maj'', ''F:maj/A'', ''C:7'']]'
$ws.Range("E15").Value = '[[''Eb:7'', ''Ab'', ''Ab''], [''Ab'', ''Ab'', ''Eb:7'']]'
$ws.Range("F15").Value = '[(17.36, 19.54), (17.94, 20.12)]'
$ws.Range("G15").Value = '[(3.1, 5.91), (7.05, 9.89)]'
$ws.Range("H15").Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'
$ws.Range("I15").Value = ""

# Row 16
$ws.Range("A16").Value = 'schubert-winterreise_170'
$ws.Range("B16").Value = 'jaah_25'
$ws.Range("C16").Value = 0.05375773651635721
$ws.Range("D16").Value = '[[''D:7'', ''G:maj'', ''D:7/C'', ''G:maj/B'']]'
$ws.Range("E16").Value = '[[''F:7'', ''Bb'', ''F:7'', ''Bb'']]'
$ws.Range("F16").Value = '[(70.38, 83.28)]'
$ws.Range("G16").Value = '[(48.25, 53.91)]'
$ws.Range("H16").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Range("I16").Value = ""

# Row 17
$ws.Range("A17").Value = 'schubert-winterreise_154'
$ws.Range("B17").Value = 'schubert-winterreise_189'
$ws.Range("C17").Value = 0.3506493506493507
$ws.Range("D17").Value = '[[''A:maj/E'', ''E:7'', ''A:maj'']]'
$ws.Range("E17").Value = '[[''G:maj/D'', ''D:7'', ''G:maj'']]'
$ws.Range("F17").Value = '[(21.78, 25.3)]'
$ws.Range("G17").Value = '[(58.32, 59.68)]'
$ws.Range("H17").Value = 'spotify:track:0XfunCHFEeQnzm4NaY8rJr'
$ws.Range("I17").Value = ""

